$d = $word.ActiveDocument

$replacements = @(
    @("268×5=1340", "270×5=1350"),
    @("794×9=7146", "749×4=2996"),
    @("621×7=4347", "824×4=3296"),
    @("250×6=1500", "986×7=6902"),
    @("679×2=1358", "883×5=4415"),
    @("624×8=4992", "625×2=1250"),
    @("529×9=4761", "596×2=1192"),
    @("466×5=2330", "780×3=2340"),
    @("336×8=2688", "890×7=6230"),
    @("536×4=2144", "980×3=2940"),
    @("291×3=873", "938×3=2814"),
    @("613×2=1226", "189×9=1701"),
    @("257×5=1285", "315×4=1260"),
    @("593×4=2372", "125×7=875"),
    @("597×4=2388", "431×3=1293"),
    @("275×5=1375", "966×5=4830"),
    @("112×6=672", "255×6=1530"),
    @("856×9=7704", "536×5=2680"),
    @("456×6=2736", "939×5=4695"),
    @("439×5=2195", "413×7=2891"),
    @("922×2=1844", "897×3=2691"),
    @("354×4=1416", "149×4=596"),
    @("736×6=4416", "912×9=8208"),
    @("582×2=1164", "889×2=1778"),
    @("150×4=600", "647×2=1294")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
